# Build dynamic lexer/parser experiment: extra rows with nested-condition
# formulas (AND / OR), a new underlined style, and tweaked pressure readings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tweak the existing "PM_IPA_FERMENTACION_PRESION" rows (2-4): the
# pressure reading in column B moves from 7.15 to 2.1 -----------------------
$ws.Range("B2").Value = 2.1
$ws.Range("B3").Value = 2.1
$ws.Range("B4").Value = 2.1

# --- Row 6: AND() nested condition -----------------------------------------
$ws.Range("A5").Copy($ws.Range("A6")) | Out-Null
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Formula = '=IF(AND(B6>3,B6<7),"presion no estable","presion estable")'
$ws.Range("E5").Copy($ws.Range("E6")) | Out-Null
$ws.Range("F5").Copy($ws.Range("F6")) | Out-Null

# --- Row 7: OR() nested condition ------------------------------------------
$ws.Range("A5").Copy($ws.Range("A7")) | Out-Null
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 1
$ws.Range("D7").Formula = '=IF(OR(B7<3,B7>7),"presion no estable","presion estable")'
$ws.Range("E5").Copy($ws.Range("E7")) | Out-Null
$ws.Range("F5").Copy($ws.Range("F7")) | Out-Null

# --- Row 11: lone styled (underlined-font) placeholder cell ----------------
$ws.Range("A11").Font.Underline = $true

# --- Selection cursor matches the author's last position -------------------
$ws.Range("A9").Select() | Out-Null
